$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.832.89'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +6.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.475.73'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +7.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +5.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.39'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +10.89%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.594'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.477.33'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.37%  '

$ws.Range("E10").Value = '  +5.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.94%  '

$ws.Range("E12").Value = '  +1.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.356'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +7.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.52'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +14.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.912.02'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.659.53'
$ws.Range("D16").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.466.05'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +6.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.36'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +8.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.87'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +10.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.35'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +8.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.84'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +5.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.89'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.176'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.55'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +15.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.29'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +7.38%  '

$ws.Range("E29").Value = '  +12.91%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0824'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +14.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.89'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +18.68%  '

$ws.Range("E32").Value = '  +9.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '175.27'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.52'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +12.36%  '

$ws.Range("E35").Value = '  +5.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.05'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +7.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '375.53'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +20.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.50'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +10.99%  '

$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("E40").Value = '  +15.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.46'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '150.74'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +10.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.75'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +9.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.86'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +12.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.600'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +6.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0969'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0528'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +8.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0229'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +7.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.25'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +9.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0231'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.29%  '
